$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper note: "Price" (column D) values are stored as TEXT in the source data
# (e.g. "26.253.23" is not a real number). Writing such numeric-looking strings
# via .Value would make Excel auto-convert them to numbers, so we force text by
# prefixing with an apostrophe, then reset .Style so no stray quote-prefix format
# lingers on the cell.

$ws.Cells.Item(2, 4).Value = "'26.253.23"
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  +0.06%  '

$ws.Cells.Item(3, 4).Value = "'1.597.16"
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  +0.54%  '

$ws.Cells.Item(4, 5).Value = '  +0.04%  '

$ws.Cells.Item(5, 4).Value = "'211.38"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.25%  '

$ws.Cells.Item(6, 4).Value = "'0.505"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +0.27%  '

$ws.Cells.Item(7, 5).Value = '  +0.06%  '

$ws.Cells.Item(8, 4).Value = "'0.245"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +0.13%  '

$ws.Cells.Item(9, 5).Value = '  +0.24%  '

$ws.Cells.Item(10, 4).Value = "'18.98"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -1.19%  '

$ws.Cells.Item(11, 5).Value = '  +0.82%  '

$ws.Cells.Item(12, 4).Value = "'1.822.08"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +0.54%  '

$ws.Cells.Item(13, 4).Value = "'1.601.31"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +0.94%  '

$ws.Cells.Item(14, 5).Value = '  -0.33%  '

$ws.Cells.Item(15, 4).Value = "'0.504"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -2.07%  '

$ws.Cells.Item(16, 4).Value = "'63.65"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -0.33%  '

$ws.Cells.Item(17, 4).Value = "'26.262.72"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +0.08%  '

$ws.Cells.Item(18, 4).Value = "'230.89"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +7.81%  '

$ws.Cells.Item(19, 5).Value = '  +3.01%  '

$ws.Cells.Item(20, 4).Value = "'0.0₃0722"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -0.48%  '

$ws.Cells.Item(21, 5).Value = '  +0.17%  '

$ws.Cells.Item(22, 4).Value = "'4.23"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.39%  '

$ws.Cells.Item(23, 4).Value = "'8.95"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -0.32%  '

$ws.Cells.Item(24, 5).Value = '  +1.31%  '

$ws.Cells.Item(25, 4).Value = "'146.27"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +1.40%  '

$ws.Cells.Item(26, 5).Value = '  +0.08%  '

$ws.Cells.Item(27, 5).Value = '  +0.32%  '

$ws.Cells.Item(28, 5).Value = '  +0.29%  '

$ws.Cells.Item(29, 4).Value = "'15.33"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +1.58%  '

$ws.Cells.Item(30, 4).Value = "'0.0493"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -0.40%  '

$ws.Cells.Item(31, 5).Value = '  -0.47%  '

$ws.Cells.Item(32, 5).Value = '  +0.61%  '

$ws.Cells.Item(33, 4).Value = "'1.471.56"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +3.91%  '

$ws.Cells.Item(34, 4).Value = "'2.95"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +0.16%  '

$ws.Cells.Item(35, 5).Value = '  -0.52%  '

$ws.Cells.Item(36, 5).Value = '  +0.46%  '

$ws.Cells.Item(37, 4).Value = "'0.568"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -3.88%  '

$ws.Cells.Item(38, 4).Value = "'0.0165"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -0.83%  '

$ws.Cells.Item(39, 4).Value = "'0.822"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -0.03%  '

$ws.Cells.Item(40, 4).Value = "'5.75"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -2.82%  '

$ws.Cells.Item(41, 5).Value = '  +0.09%  '

$ws.Cells.Item(42, 5).Value = '  +2.03%  '

$ws.Cells.Item(43, 4).Value = "'0.930"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -1.11%  '

$ws.Cells.Item(44, 4).Value = "'1.734.64"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +0.61%  '

$ws.Cells.Item(45, 4).Value = "'0.757"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -1.16%  '

$ws.Cells.Item(46, 4).Value = "'60.63"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -0.85%  '

$ws.Cells.Item(47, 4).Value = "'87.98"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +2.56%  '

$ws.Cells.Item(48, 5).Value = '  -0.74%  '

$ws.Cells.Item(49, 4).Value = "'0.0502"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +0.05%  '

$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(50, 4).Value = "'0.0949"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -1.90%  '

$ws.Cells.Item(51, 2).Value = 'USDD'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Cells.Item(51, 4).Value = "'0.998"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  -0.06%  '
